# Remove the trailing "Ver no Jupiter ..." and "© 2020 ..." paragraphs
# (plus the blank paragraph that separated them from the bibliography),
# leaving the bibliography's last line, one blank paragraph, and the
# final page-break paragraph intact.

$d = $word.ActiveDocument

# Anchor on the last bibliography line that must be kept.
$anchorEnd = $d.Content
$foundStart = $anchorEnd.Find.Execute(
    "uma introdução à álgebra linear. São Paulo: Thomson, 2007.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Anchor on the copyright line, the last paragraph that must be removed.
$anchorFooter = $d.Content
$foundEnd = $anchorFooter.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart -and $foundEnd) {
    # Start right after the bibliography line's paragraph mark (i.e. at
    # the start of the blank paragraph that follows it) and end right
    # after the copyright paragraph's own paragraph mark, so the whole
    # blank/"Ver no Jupiter.../"© 2020..." block (3 paragraphs) is removed
    # as one unit, merging back into the blank paragraph that follows.
    $deleteStart = $anchorEnd.End + 1
    $deleteEnd = $anchorFooter.End + 1

    $victim = $d.Range($deleteStart, $deleteEnd)
    $victim.Delete()
}
